# Update "想去人数" (F column) counts on the 展览 sheet and the 全部类型 sheet
# to reflect output generated at a later point in time.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 2242
$ws1.Range("F9").Value  = 291
$ws1.Range("F21").Value = 4104
$ws1.Range("F23").Value = 737
$ws1.Range("F24").Value = 34
$ws1.Range("F36").Value = 2661

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 2242
$ws4.Range("F9").Value  = 291
$ws4.Range("F21").Value = 4104
$ws4.Range("F23").Value = 737
$ws4.Range("F24").Value = 34
$ws4.Range("F37").Value = 2661
